# GSC export refresh: drop the incomplete leading day ("2025-10-15", which
# had no data yet) from the "Chart" data table. Deleting the row shifts
# every subsequent date up by one, so the table now runs from 2025-10-16
# through 2026-01-12 (previously the table ran one extra, trailing day
# past 2026-01-12 that carried a stray blank "Impressions" note).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (row 2 = 2025-10-15); Excel shifts rows 3..91
# up into rows 2..90.
$ws.Rows.Item(2).Delete()

# The last row (now row 90, date 2026-01-12) previously carried the next
# row's trailing blank-string "Impressions" note; normalize it back to a
# plain numeric 0, matching every other row's Impressions column.
$ws.Range("D90").Value = 0
